$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.532.73"
$ws.Range("E2").Value = "  +1.85%  "

$ws.Range("D3").Value = "2.290.01"
$ws.Range("E3").Value = "  +1.22%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'156.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +15,556.28%  "

$ws.Range("D6").Value = "'307.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.03%  "

$ws.Range("D7").Value = "'96.12"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.06%  "

$ws.Range("E8").Value = "  +0.17%  "

$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("D10").Value = "'0.497"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.31%  "

$ws.Range("D11").Value = "'35.83"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +10.76%  "

$ws.Range("D12").Value = "'0.0807"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.66%  "

$ws.Range("E13").Value = "  -1.61%  "

$ws.Range("D14").Value = "'6.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.49%  "

$ws.Range("D15").Value = "2.643.95"
$ws.Range("E15").Value = "  +1.20%  "

$ws.Range("E16").Value = "  +2.59%  "

$ws.Range("D17").Value = "2.288.62"
$ws.Range("E17").Value = "  +0.26%  "

$ws.Range("D18").Value = "'0.804"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.82%  "

$ws.Range("D19").Value = "42.378.70"
$ws.Range("E19").Value = "  +1.72%  "

$ws.Range("D20").Value = "'12.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.53%  "

$ws.Range("D21").Value = "0.0₃0924"
$ws.Range("E21").Value = "  +2.26%  "

$ws.Range("E22").Value = "  +2.28%  "

$ws.Range("D23").Value = "'68.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.29%  "

$ws.Range("D24").Value = "'243.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.30%  "

$ws.Range("D25").Value = "'2.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.90%  "

$ws.Range("E26").Value = "  +2.58%  "

$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.21%  "

$ws.Range("D28").Value = "'24.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.51%  "

$ws.Range("D29").Value = "'36.14"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.03%  "

$ws.Range("D30").Value = "'9.66"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.22%  "

$ws.Range("E31").Value = "  -8.41%  "

$ws.Range("D32").Value = "'161.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.98%  "

$ws.Range("D33").Value = "'5.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.52%  "

$ws.Range("E34").Value = "  -0.06%  "

$ws.Range("E35").Value = "  +1.69%  "

$ws.Range("E36").Value = "  +3.19%  "

$ws.Range("E37").Value = "  +5.38%  "

$ws.Range("D38").Value = "'17.30"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.26%  "

$ws.Range("E39").Value = "  +0.33%  "

$ws.Range("E40").Value = "  +0.06%  "

$ws.Range("D41").Value = "'1.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.12%  "

$ws.Range("D42").Value = "'4.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.09%  "

$ws.Range("D43").Value = "2.021.77"
$ws.Range("E43").Value = "  -1.91%  "

$ws.Range("D44").Value = "'19.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.63%  "

$ws.Range("E45").Value = "  +11.09%  "

$ws.Range("E46").Value = "  +2.44%  "

$ws.Range("E47").Value = "  +0.30%  "

$ws.Range("D48").Value = "'3.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.60%  "

$ws.Range("E49").Value = "  +1.91%  "

$ws.Range("D50").Value = "'53.47"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.97%  "

$ws.Range("D51").Value = "'73.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.20%  "
